$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# Title / heading (appears twice: Heading1 and the bold summary line)
Replace-Text "Play Ice Wolf Free Slot Game Review | Enjoy Unique Game Features" "Play Ice Wolf Free Online | Slot Game Review"

# "What we like" bullet points
Replace-Text "Unique and engaging game mechanics" "Unique gameplay mechanics and features"
Replace-Text "Stunning graphics and theme design" "Stunning setting and theme design"
Replace-Text "Special symbols and animations add excitement" "Engaging game symbol and animation"
Replace-Text "Fantastic background music enhances the experience" "Immersive sound and music design"

# "What we don't like" bullet points
Replace-Text "Fewer winnings during regular spins" "Reduced number of winnings during spins"
Replace-Text "Top ice must be unlocked to earn high-value winnings" "Limited selection of similar themed slots"

# Meta description (italic summary paragraph)
Replace-Text "Read our review for Ice Wolf, the online slot game with engaging mechanics and stunning graphics. Play for free and experience special symbols and animations." "Read our review of Ice Wolf, an immersive slot game. Play for free and enjoy unique gameplay mechanics and stunning design."
